# Generate Report for Handoff
#
# This re-generates the localization-status report's timestamp columns:
# the "Latest Handoff Date"/"Latest Handoff Datetime" values that were
# last refreshed at 2016-03-23 06:30:11 / 06:30:51 are refreshed to
# 2016-03-23 06:31:50, and the values last refreshed at
# 2016-03-23 06:30:01 / 06:30:44 are refreshed to 2016-03-23 06:31:44.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D ("Latest Handoff Date") ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value  = "2016-03-23 06:31:50"
$overview.Range("D10").Value = "2016-03-23 06:31:50"
$overview.Range("D11").Value = "2016-03-23 06:31:50"
$overview.Range("D12").Value = "2016-03-23 06:31:50"
$overview.Range("D13").Value = "2016-03-23 06:31:50"
$overview.Range("D14").Value = "2016-03-23 06:31:50"
$overview.Range("D15").Value = "2016-03-23 06:31:50"
$overview.Range("D16").Value = "2016-03-23 06:31:50"

# --- zh-cn sheet: column E ("Latest Handoff Datetime") ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value  = "2016-03-23 06:31:44"
$zhcn.Range("E10").Value = "2016-03-23 06:31:44"
$zhcn.Range("E11").Value = "2016-03-23 06:31:44"
$zhcn.Range("E12").Value = "2016-03-23 06:31:44"
$zhcn.Range("E13").Value = "2016-03-23 06:31:44"
$zhcn.Range("E14").Value = "2016-03-23 06:31:44"
$zhcn.Range("E15").Value = "2016-03-23 06:31:44"
$zhcn.Range("E16").Value = "2016-03-23 06:31:44"

# --- de-de sheet: column E ("Latest Handoff Datetime") ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E7").Value  = "2016-03-23 06:31:50"
$dede.Range("E10").Value = "2016-03-23 06:31:50"
$dede.Range("E11").Value = "2016-03-23 06:31:50"
$dede.Range("E12").Value = "2016-03-23 06:31:50"
$dede.Range("E13").Value = "2016-03-23 06:31:50"
$dede.Range("E14").Value = "2016-03-23 06:31:50"
$dede.Range("E15").Value = "2016-03-23 06:31:50"
$dede.Range("E16").Value = "2016-03-23 06:31:50"
